$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared "short-url" value (column B) used by every data row (2..558)
# changes from "Xpxm4b" to "9UzbA5".
$ws.Range("B2:B558").Value = "9UzbA5"

# Row 540: refugees (N) 540 -> 505 ; asylum_seekers (O) 240 -> 267
$ws.Range("N540").Value = 505
$ws.Range("O540").Value = 267

# Row 541: asylum_seekers (O) 5 -> 6
$ws.Range("O541").Value = 6

# Row 543: refugees (N) 19 -> 20
$ws.Range("N543").Value = 20

# Row 544: refugees (N) 8258 -> 8045 ; asylum_seekers (O) 3766 -> 4102 ; ooc (T) 668 -> 666
$ws.Range("N544").Value = 8045
$ws.Range("O544").Value = 4102
$ws.Range("T544").Value = 666

# Row 545: refugees (N) 7 -> 6
$ws.Range("N545").Value = 6

# Row 546: refugees (N) 72 -> 61 ; asylum_seekers (O) 24 -> 31
$ws.Range("N546").Value = 61
$ws.Range("O546").Value = 31

# Row 548: asylum_seekers (O) 123 -> 124
$ws.Range("O548").Value = 124

# Row 549: asylum_seekers (O) 8405 -> 8455
$ws.Range("O549").Value = 8455

# Row 551: refugees (N) 576 -> 535 ; asylum_seekers (O) 84 -> 90
$ws.Range("N551").Value = 535
$ws.Range("O551").Value = 90

# Row 553: asylum_seekers (O) 20 -> 21
$ws.Range("O553").Value = 21

# Row 558: ooc (T) 128 -> 126
$ws.Range("T558").Value = 126
